$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - copy formatting (font/border/alignment) from E1, then set the value
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 10:51:57.334322",
    "2021-10-05 10:51:57.334334",
    "2021-10-05 10:51:57.334338",
    "2021-10-05 10:51:57.334341",
    "2021-10-05 10:51:57.334345",
    "2021-10-05 10:51:57.334348",
    "2021-10-05 10:51:57.334351",
    "2021-10-05 10:51:57.334354",
    "2021-10-05 10:51:57.334357",
    "2021-10-05 10:51:57.334360",
    "2021-10-05 10:51:57.334363",
    "2021-10-05 10:51:57.334366",
    "2021-10-05 10:51:57.334369",
    "2021-10-05 10:51:57.334372",
    "2021-10-05 10:51:57.334375",
    "2021-10-05 10:51:57.334378",
    "2021-10-05 10:51:57.334381",
    "2021-10-05 10:51:57.334384",
    "2021-10-05 10:51:57.334402",
    "2021-10-05 10:51:57.334413",
    "2021-10-05 10:51:57.334418",
    "2021-10-05 10:51:57.334436",
    "2021-10-05 10:51:57.334475",
    "2021-10-05 10:51:57.334479",
    "2021-10-05 10:51:57.334483",
    "2021-10-05 10:51:57.334487",
    "2021-10-05 10:51:57.334490",
    "2021-10-05 10:51:57.334493",
    "2021-10-05 10:51:57.334496",
    "2021-10-05 10:51:57.334499",
    "2021-10-05 10:51:57.334502",
    "2021-10-05 10:51:57.334505",
    "2021-10-05 10:51:57.334509",
    "2021-10-05 10:51:57.334513",
    "2021-10-05 10:51:57.334516",
    "2021-10-05 10:51:57.334519",
    "2021-10-05 10:51:57.334522",
    "2021-10-05 10:51:57.334525",
    "2021-10-05 10:51:57.334528",
    "2021-10-05 10:51:57.334531",
    "2021-10-05 10:51:57.334535",
    "2021-10-05 10:51:57.334538",
    "2021-10-05 10:51:57.334541",
    "2021-10-05 10:51:57.334544",
    "2021-10-05 10:51:57.334547",
    "2021-10-05 10:51:57.334550",
    "2021-10-05 10:51:57.334553",
    "2021-10-05 10:51:57.334556",
    "2021-10-05 10:51:57.334559",
    "2021-10-05 10:51:57.334562",
    "2021-10-05 10:51:57.334565",
    "2021-10-05 10:51:57.334568",
    "2021-10-05 10:51:57.334572",
    "2021-10-05 10:51:57.334575",
    "2021-10-05 10:51:57.334578"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

